$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue "D2" "29.148.47"
Set-TextValue "E2" "  -3.29%  "
Set-TextValue "D3" "1.850.70"
Set-TextValue "E3" "  -2.23%  "
Set-TextValue "D4" "0.9998"
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "0.7024"
Set-TextValue "E5" "  -5.25%  "
Set-TextValue "D6" "238.30"
Set-TextValue "E7" "  +0.03%  "
Set-TextValue "D8" "0.3042"
Set-TextValue "E8" "  -4.23%  "
Set-TextValue "D9" "0.07513"
Set-TextValue "E9" "  +4.00%  "
Set-TextValue "D10" "23.38"
Set-TextValue "E10" "  -6.42%  "
Set-TextValue "D11" "0.08131"
Set-TextValue "E11" "  -2.57%  "
Set-TextValue "D12" "1.852.66"
Set-TextValue "E12" "  -4.47%  "
Set-TextValue "D13" "0.7249"
Set-TextValue "E13" "  -4.78%  "
Set-TextValue "D14" "5.219"
Set-TextValue "E14" "  -4.31%  "
Set-TextValue "D15" "88.64"
Set-TextValue "E15" "  -5.01%  "
Set-TextValue "D16" "29.267.92"
Set-TextValue "E16" "  -3.11%  "
Set-TextValue "D17" "5.763"
Set-TextValue "D18" "237.40"
Set-TextValue "E18" "  -5.39%  "
Set-TextValue "D19" "13.07"
Set-TextValue "E19" "  -4.25%  "
Set-TextValue "D20" "0.000007632"
Set-TextValue "E20" "  -3.44%  "
Set-TextValue "D21" "1.000"
Set-TextValue "E21" "  -0.16%  "
Set-TextValue "D22" "2.118.37"
Set-TextValue "E22" "  -4.64%  "
Set-TextValue "E23" "  +0.09%  "
Set-TextValue "D24" "7.560"
Set-TextValue "E24" "  -5.20%  "
Set-TextValue "D25" "8.988"
Set-TextValue "D26" "161.19"
Set-TextValue "E26" "  -2.10%  "
Set-TextValue "D27" "0.1453"
Set-TextValue "E27" "  -8.27%  "
Set-TextValue "D28" "18.05"
Set-TextValue "E28" "  -3.97%  "
Set-TextValue "D29" "1.968"
Set-TextValue "E29" "  -4.98%  "
Set-TextValue "D30" "1.399"
Set-TextValue "E30" "  -6.01%  "
Set-TextValue "D31" "4.522"
Set-TextValue "E31" "  -1.62%  "
Set-TextValue "D32" "1.493"
Set-TextValue "E32" "  -3.07%  "
Set-TextValue "D33" "3.970"
Set-TextValue "E33" "  -5.71%  "
Set-TextValue "D34" "0.05153"
Set-TextValue "E34" "  -4.23%  "
Set-TextValue "D35" "1.185"
Set-TextValue "B36" "ImmutableX"
Set-TextValue "C36" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D36" "0.6994"
Set-TextValue "E36" "  -10.11%  "
Set-TextValue "B37" "Frax"
Set-TextValue "C37" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D37" "1.010"
Set-TextValue "E37" "  +0.46%  "
Set-TextValue "D38" "2.657"
Set-TextValue "D39" "0.01861"
Set-TextValue "E39" "  -5.25%  "
Set-TextValue "E40" "  -3.18%  "
Set-TextValue "D41" "0.9374"
Set-TextValue "E41" "  +5.64%  "
Set-TextValue "E42" "  -1.97%  "
Set-TextValue "D43" "1.076.98"
Set-TextValue "E43" "  -2.25%  "
Set-TextValue "D44" "0.4281"
Set-TextValue "E44" "  -6.50%  "
Set-TextValue "D45" "69.87"
Set-TextValue "E45" "  -3.85%  "
Set-TextValue "D46" "0.9997"
Set-TextValue "E46" "  -0.14%  "
Set-TextValue "D47" "102.24"
Set-TextValue "E47" "  -2.22%  "
Set-TextValue "D48" "1.742"
Set-TextValue "E48" "  -6.80%  "
Set-TextValue "D49" "1.991.37"
Set-TextValue "E49" "  -5.72%  "
Set-TextValue "D50" "7.041"
Set-TextValue "E50" "  -7.45%  "
Set-TextValue "D51" "9.118"
Set-TextValue "E51" "  -5.48%  "
